$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 11:58"

# Country name corrections - swap mislabeled rows so the country list
# matches the correct case-count ordering
$ws.Range("A21").Value = "Indonesia"
$ws.Range("A22").Value = "Filipinas"
$ws.Range("A34").Value = "Polonia"
$ws.Range("A35").Value = "Ecuador"
$ws.Range("A53").Value = "Bielorrusia"
$ws.Range("A54").Value = "Honduras"

# Updated COVID-19 statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Muertes hoy, Muertes)
$ws.Range("B4").Value = 8219088
$ws.Range("C4").Value = 2773
$ws.Range("E4").Value = 2675948
$ws.Range("G4").Value = 37
$ws.Range("H4").Value = 222754
$ws.Range("B19").Value = 386086
$ws.Range("C19").Value = 1527
$ws.Range("D19").Value = 300738
$ws.Range("E19").Value = 79725
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = 5623
$ws.Range("B21").Value = 353461
$ws.Range("C21").Value = 4301
$ws.Range("D21").Value = 277544
$ws.Range("E21").Value = 63570
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = 12347
$ws.Range("B22").Value = 351750
$ws.Range("C22").Value = 3139
$ws.Range("D22").Value = 294865
$ws.Range("E22").Value = 50354
$ws.Range("G22").Value = 34
$ws.Range("H22").Value = 6531
$ws.Range("B23").Value = 349639
$ws.Range("C23").Value = 823
$ws.Range("E23").Value = 55223
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 9816
$ws.Range("B27").Value = 301389
$ws.Range("C27").Value = 1188
$ws.Range("D27").Value = 261128
$ws.Range("E27").Value = 38122
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = 2139
$ws.Range("B34").Value = 157608
$ws.Range("C34").Value = 7705
$ws.Range("D34").Value = 87773
$ws.Range("E34").Value = 66395
$ws.Range("G34").Value = 132
$ws.Range("H34").Value = 3440
$ws.Range("B35").Value = 150360
$ws.Range("D35").Value = 128134
$ws.Range("E35").Value = 9920
$ws.Range("H35").Value = 12306
$ws.Range("B53").Value = 86392
$ws.Range("C53").Value = 658
$ws.Range("D53").Value = 78990
$ws.Range("E53").Value = 6481
$ws.Range("G53").Value = 5
$ws.Range("H53").Value = 921
$ws.Range("B54").Value = 86089
$ws.Range("C54").Value = 631
$ws.Range("D54").Value = 34099
$ws.Range("E54").Value = 49438
$ws.Range("G54").Value = 19
$ws.Range("H54").Value = 2552
$ws.Range("B83").Value = 31265
$ws.Range("C83").Value = 204
$ws.Range("D83").Value = 26542
$ws.Range("E83").Value = 3811
$ws.Range("B102").Value = 13133
$ws.Range("C102").Value = 189
$ws.Range("E102").Value = 3683
$ws.Range("B104").Value = 11327
$ws.Range("C104").Value = 72
$ws.Range("D104").Value = 10380
$ws.Range("E104").Value = 877
$ws.Range("D128").Value = 3385
$ws.Range("E128").Value = 1846
$ws.Range("B129").Value = 5221
$ws.Range("C129").Value = 7
$ws.Range("D129").Value = 4951
$ws.Range("E129").Value = 165
$ws.Range("B195").Value = 149
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 148
$ws.Range("E195").Value = 1
